$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell while preventing Excel from
# auto-coercing numeric-looking text (e.g. "1.210" or "3.686") into a
# floating point number, which would silently drop meaningful trailing
# zeros / precision. Values containing two or more '.' characters (e.g.
# "28.083.39") are never parsed as numbers by Excel, so they are left
# with the default (General) format.
function Set-TextSafeValue {
    param($cell, [string]$val)
    $dotCount = ($val.ToCharArray() | Where-Object { $_ -eq '.' }).Count
    if ($dotCount -le 1) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

$updates = @(
    @{Row=2; D='28.083.39'; E='  -1.34%  '},
    @{Row=3; D='1.789.10'; E='  -1.69%  '},
    @{Row=4; D='1.004'; E='  +0.26%  '},
    @{Row=5; D='313.51'; E='  -0.71%  '},
    @{Row=6; D='1.003'; E='  +0.21%  '},
    @{Row=7; D='0.5197'; E='  +2.10%  '},
    @{Row=8; D='0.3797'; E='  -3.91%  '},
    @{Row=9; D='0.07962'; E='  -3.83%  '},
    @{Row=10; D='41.37'; E='  -0.54%  '},
    @{Row=11; D='1.089'; E='  -1.99%  '},
    @{Row=12; D='1.004'; E='  +0.27%  '},
    @{Row=13; D='6.251'; E='  -1.01%  '},
    @{Row=14; D='20.46'; E='  -2.77%  '},
    @{Row=15; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.791.35'; E='  -1.28%  '},
    @{Row=16; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.270'; E='  -3.35%  '},
    @{Row=17; D='91.14'; E='  -1.56%  '},
    @{Row=18; D='0.00001089'; E='  -5.11%  '},
    @{Row=19; D='0.06541'; E='  -1.63%  '},
    @{Row=20; D='1.003'; E='  +0.30%  '},
    @{Row=21; D='17.26'; E='  -2.82%  '},
    @{Row=22; E='  -2.97%  '},
    @{Row=23; D='28.135.50'; E='  -1.22%  '},
    @{Row=24; D='11.09'; E='  -3.19%  '},
    @{Row=25; D='2.260'; E='  -0.52%  '},
    @{Row=26; D='159.71'; E='  +2.42%  '},
    @{Row=27; D='20.34'; E='  -4.55%  '},
    @{Row=28; D='1.993.01'; E='  -1.55%  '},
    @{Row=29; D='2.323'; E='  -3.61%  '},
    @{Row=30; E='  -2.65%  '},
    @{Row=31; D='0.1083'; E='  -1.05%  '},
    @{Row=32; D='1.052'; E='  -5.34%  '},
    @{Row=33; D='3.686'},
    @{Row=34; D='5.521'; E='  -4.40%  '},
    @{Row=35; D='0.07176'; E='  +1.73%  '},
    @{Row=36; D='12.13'; E='  +7.51%  '},
    @{Row=37; D='0.02303'; E='  -1.55%  '},
    @{Row=38; D='0.2134'; E='  -4.15%  '},
    @{Row=39; D='5.062'; E='  -3.62%  '},
    @{Row=40; D='8.557'; E='  -3.57%  '},
    @{Row=41; D='0.6144'; E='  -2.53%  '},
    @{Row=42; E='  -1.53%  '},
    @{Row=43; D='1.370'; E='  -2.27%  '},
    @{Row=44; E='  -1.99%  '},
    @{Row=45; D='3.756'; E='  +0.59%  '},
    @{Row=46; D='0.5904'; E='  -0.39%  '},
    @{Row=47; D='126.61'; E='  +1.11%  '},
    @{Row=48; D='1.210'; E='  +2.05%  '},
    @{Row=49; D='1.913'; E='  -3.56%  '},
    @{Row=50; D='0.06759'; E='  -1.93%  '},
    @{Row=51; D='72.30'; E='  -2.70%  '}
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.ContainsKey('B')) {
        $ws.Range("B$row").Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Range("C$row").Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        Set-TextSafeValue $ws.Range("D$row") $item.D
    }
    if ($item.ContainsKey('E')) {
        $ws.Range("E$row").Value = $item.E
    }
}
